# Update the "Through ..." sheet name to reflect the new data date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-06-09"

# Update the header label for the June column (shared string)
$ws.Range("A7").Value = "June (through 06-09)"

# Update the June row (row 7) values for columns C through I (2016-2022)
$ws.Range("C7").Value = 13
$ws.Range("D7").Value = 19
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 41
$ws.Range("I7").Value = 34

# Update the Total row (row 8) values for columns C through I (2016-2022)
$ws.Range("C8").Value = 222
$ws.Range("D8").Value = 335
$ws.Range("E8").Value = 314
$ws.Range("F8").Value = 214
$ws.Range("G8").Value = 399
$ws.Range("I8").Value = 697
